$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "Shared Formula"

# E2: plain (non-shared) formula referencing B2
$ws.Range("E2").Formula = "=B2 + 1"

# E3:E4 use the same relative formula pattern, entered together so Excel
# records them as a shared formula group (t="shared").
$ws.Range("E3:E4").Formula = "=B3 + 1"

# Match the column width Excel computed for the new column (bestFit-style
# autosize to fit the "Shared Formula" header)
$ws.Columns.Item(5).ColumnWidth = 12.25

# Leave the selection where the author left off
$ws.Range("E3").Select()
